$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62: fill in the Stop time, Delta minutes, and Activity Notes that were
# left blank, completing the "Debugging" log entry.
$ws.Range("C62").Value = 0.15277777777777776
$ws.Range("C62").NumberFormat = "h:mm"
$ws.Range("E62").Value = 47
$ws.Range("H62").Value = "Debugged inserts into Orders table in db. All good now"

# New row 63: a fresh log entry for the next coding session.
$ws.Range("B63").Value = "3:40PM"
$ws.Range("F63").Value = "Code "
$ws.Range("G63").Value = "Adding Products to OrderItems table from given products in order"

# Update the sheet view: scroll position and active selection move to the
# newly added row.
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 2
$ws.Range("G63").Select() | Out-Null
